# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E23) is re-sorted from descending
# (2005 .. 1910) to ascending (1910 .. 2005), and the "Valor Mora"
# column (F16:F23) is kept in sync with its period row (the values for
# 1910 and 2005 swap places as a result).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1910", "1911", "1912", "2001", "2002", "2003", "2004", "2005")
$values  = @(14400, 48000, 48000, 48000, 48000, 48000, 48000, 30400)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
